$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 281
$ws1.Range("F3").Value = 185
$ws1.Range("F4").Value = 2152
$ws1.Range("F5").Value = 1673
$ws1.Range("F6").Value = 309
$ws1.Range("F7").Value = 93
$ws1.Range("F8").Value = 733
$ws1.Range("F9").Value = 160

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 22

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 281
$ws4.Range("F3").Value = 185
$ws4.Range("F4").Value = 2152
$ws4.Range("F5").Value = 1673
$ws4.Range("F6").Value = 309
$ws4.Range("F7").Value = 22
$ws4.Range("F8").Value = 93
$ws4.Range("F9").Value = 733
$ws4.Range("F10").Value = 160
